$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New I/O entry: housePowerLink / data / comms link note
$ws.Range("B23").Value = "housePowerLink"
$ws.Range("D23").Value = "data"
$ws.Range("E23").Value = "comms link with the house power battery monitor"

# D23's "Type" column should be centered like the other "data" rows (D20:D22).
# Touching NumberFormat first keeps the alignment write from pulling in a
# stray border from an unrelated existing style.
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").HorizontalAlignment = -4108

# Move the active selection to B24, matching the saved view state
$ws.Range("B24").Select()

$wb.Save()
